# "Added last minute updates"
#
# The document's opening paragraph holds a merge-field style placeholder
# (**ID__AFFARS_..._ID**) followed by a run that is just a single space.
# This edit:
#   1. Repoints the placeholder ID to the new topic id.
#   2. Drops the now unnecessary trailing-space run.
#   3. Adds a 5-twip paragraph border (top/left/bottom/right, line-less —
#      i.e. only the `w:space` padding is set) around that paragraph.
#   4. Bumps the paragraph's left indent from 120 to 225 twips.

$d = $word.ActiveDocument

# --- 1. Swap the placeholder ID text -----------------------------------
$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5315_topic_3__ID**", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "**ID__AFFARS_AF_PGI_5315_209_90__ID**", 2
) | Out-Null

# --- 2. Remove the trailing " " run from the first paragraph -----------
$p1 = $d.Paragraphs(1)
$pEnd = $p1.Range.End
# the paragraph mark sits at [pEnd-1, pEnd); the space run is right before it
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# --- 3 & 4. Paragraph border + new left indent --------------------------
$pf = $p1.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25   # 225 twips (twentieths of a point)
